$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected answer texts (fixes the "',target=" -> "' target=" typo in the
# existing Meat/Dairy answers, and turns the Farm Produce answer into the same
# linked-list style as Meat/Dairy) -----------------------------------------
$meatAnswer = 'We have multiple kinds of Meat. Please select your choice. 1.<a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Beef</a> 2. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Pork</a> 3. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Chicken</a>'

$dairyAnswer = 'Among Dairy products we have 1. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Milk</a> 2. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Yogurt</a> 3. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Greek Yogurt</a> 4. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Cheese</a> 5. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Milk Powder</a> 6. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Butter</a> 7. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Ice Cream</a>'

$produceAnswer = 'Among Farm Produce we have 1. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Carrots</a> 2. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Tomato</a> 3. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Potato</a> 4. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Onions</a> 5. <a href = ''https://www.sysco.com/Products/Products/Product-Categories.html'' target=''_blank''>Broccili</a>'

# --- Insert a new Question/Answer row before the old row 28 ("Cool Produce")
# so a "farm produce" question is added, shifting the remaining rows (old
# 28-40, "Cool Produce" through "Crompton Executive") down to 29-41. -------
$ws.Rows(28).Insert()
$ws.Range("A28").Value = "farm produce"
$ws.Range("B28").Value = $produceAnswer

# --- Refresh the (now corrected) answer text everywhere it is used. -------
for ($r = 2; $r -le 10; $r++) {
    $ws.Range("B" + $r).Value = $meatAnswer
}
for ($r = 11; $r -le 19; $r++) {
    $ws.Range("B" + $r).Value = $dairyAnswer
}
for ($r = 20; $r -le 29; $r++) {
    $ws.Range("B" + $r).Value = $produceAnswer
}

# --- Match the saved selection state recorded in the workbook. ------------
$ws.Range("B2").Select() | Out-Null
